$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: round coordinates, drop the time cells (Z2, AB2) ---
$ws.Range("Q2").Value = 507350
$ws.Range("R2").Value = 6946859
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3: round coordinates, drop the time cells (Z3, AB3) ---
$ws.Range("Q3").Value = 507339
$ws.Range("R3").Value = 6946917
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# --- Row 4: replaced by what used to be row 5's record, coords rounded ---
$ws.Range("A4").Value = 111801760
$ws.Range("B4").Value = 89558
$ws.Range("E4").Value = 1503
$ws.Range("F4").Value = "Gräddporing"
$ws.Range("G4").Value = "Sidera lenis"
$ws.Range("H4").Value = "(P.Karst.) Miettinen"
$ws.Range("Q4").Value = 507293
$ws.Range("R4").Value = 6946996
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

# --- Row 5: replaced by what used to be row 4's record, coords rounded ---
$ws.Range("A5").Value = 111801750
$ws.Range("B5").Value = 89745
$ws.Range("E5").Value = 2062
$ws.Range("F5").Value = "Ulltickeporing"
$ws.Range("G5").Value = "Skeletocutis brevispora"
$ws.Range("H5").Value = "Niemelä"
$ws.Range("Q5").Value = 507350
$ws.Range("R5").Value = 6946859
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()

Write-Host "edit complete"
